# Update NATMI LR-pair sheet (Gdf7-Bmpr1a) with newly recomputed TPM-based
# statistics and add "Resolving-Mac" as an additional sending cluster.
#
# Rows 2-5 keep their Sending/Ligand/Receptor/Target cluster labels
# (MuSCs / Gdf7 / Bmpr1a / {ECs,FAPs,MuSCs,Resolving-Mac}); only the
# numeric expression / specificity columns (G,H,I,J,M,N,O,P,Q,R,S,T) are
# refreshed. Rows 6-9 are brand new: "Resolving-Mac" acting as the sending
# cluster against the same four target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: MuSCs | Gdf7 | Bmpr1a | ECs
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 7).Value  = 0.02913733333333333
$ws.Cells.Item(2, 8).Value  = 0.087412
$ws.Cells.Item(2, 9).Value  = 0.2991461473965196
$ws.Cells.Item(2, 10).Value = 0.2991461473965196
$ws.Cells.Item(2, 13).Value = 6.322177333333333
$ws.Cells.Item(2, 14).Value = 18.966532
$ws.Cells.Item(2, 15).Value = 0.08271011762055308
$ws.Cells.Item(2, 16).Value = 0.08271011762055309
$ws.Cells.Item(2, 17).Value = 0.1842113883537778
$ws.Cells.Item(2, 18).Value = 1.657902495184
$ws.Cells.Item(2, 19).Value = 0.02474241303690144
$ws.Cells.Item(2, 20).Value = 0.02474241303690145

# ---------------------------------------------------------------------
# Row 3: MuSCs | Gdf7 | Bmpr1a | FAPs
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 7).Value  = 0.02913733333333333
$ws.Cells.Item(3, 8).Value  = 0.087412
$ws.Cells.Item(3, 9).Value  = 0.2991461473965196
$ws.Cells.Item(3, 10).Value = 0.2991461473965196
$ws.Cells.Item(3, 13).Value = 41.286995
$ws.Cells.Item(3, 14).Value = 123.860985
$ws.Cells.Item(3, 15).Value = 0.5401386314560596
$ws.Cells.Item(3, 16).Value = 0.5401386314560597
$ws.Cells.Item(3, 17).Value = 1.202992935646667
$ws.Cells.Item(3, 18).Value = 10.82693642082
$ws.Cells.Item(3, 19).Value = 0.1615803906601088
$ws.Cells.Item(3, 20).Value = 0.1615803906601088

# ---------------------------------------------------------------------
# Row 4: MuSCs | Gdf7 | Bmpr1a | MuSCs
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 7).Value  = 0.02913733333333333
$ws.Cells.Item(4, 8).Value  = 0.087412
$ws.Cells.Item(4, 9).Value  = 0.2991461473965196
$ws.Cells.Item(4, 10).Value = 0.2991461473965196
$ws.Cells.Item(4, 13).Value = 27.73243066666667
$ws.Cells.Item(4, 14).Value = 83.197292
$ws.Cells.Item(4, 15).Value = 0.3628105447549136
$ws.Cells.Item(4, 16).Value = 0.3628105447549136
$ws.Cells.Item(4, 17).Value = 0.8080490764782223
$ws.Cells.Item(4, 18).Value = 7.272441688304001
$ws.Cells.Item(4, 19).Value = 0.1085333766982649
$ws.Cells.Item(4, 20).Value = 0.1085333766982649

# ---------------------------------------------------------------------
# Row 5: MuSCs | Gdf7 | Bmpr1a | Resolving-Mac
# ---------------------------------------------------------------------
$ws.Cells.Item(5, 7).Value  = 0.02913733333333333
$ws.Cells.Item(5, 8).Value  = 0.087412
$ws.Cells.Item(5, 9).Value  = 0.2991461473965196
$ws.Cells.Item(5, 10).Value = 0.2991461473965196
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.096171666666667
$ws.Cells.Item(5, 14).Value = 3.288515
$ws.Cells.Item(5, 15).Value = 0.01434070616847367
$ws.Cells.Item(5, 16).Value = 0.01434070616847367
$ws.Cells.Item(5, 17).Value = 0.03193951924222223
$ws.Cells.Item(5, 18).Value = 0.28745567318
$ws.Cells.Item(5, 19).Value = 0.004289967001244401
$ws.Cells.Item(5, 20).Value = 0.004289967001244401

# ---------------------------------------------------------------------
# Row 6 (new): Resolving-Mac | Gdf7 | Bmpr1a | ECs
# ---------------------------------------------------------------------
$ws.Cells.Item(6, 1).Value  = "Resolving-Mac"
$ws.Cells.Item(6, 2).Value  = "Gdf7"
$ws.Cells.Item(6, 3).Value  = "Bmpr1a"
$ws.Cells.Item(6, 4).Value  = "ECs"
$ws.Cells.Item(6, 5).Value  = 1
$ws.Cells.Item(6, 6).Value  = 0.3333333333333333
$ws.Cells.Item(6, 7).Value  = 0.06826433333333333
$ws.Cells.Item(6, 8).Value  = 0.204793
$ws.Cells.Item(6, 9).Value  = 0.7008538526034804
$ws.Cells.Item(6, 10).Value = 0.7008538526034804
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 6.322177333333333
$ws.Cells.Item(6, 14).Value = 18.966532
$ws.Cells.Item(6, 15).Value = 0.08271011762055308
$ws.Cells.Item(6, 16).Value = 0.08271011762055309
$ws.Cells.Item(6, 17).Value = 0.4315792208751111
$ws.Cells.Item(6, 18).Value = 3.884212987876
$ws.Cells.Item(6, 19).Value = 0.05796770458365164
$ws.Cells.Item(6, 20).Value = 0.05796770458365165

# ---------------------------------------------------------------------
# Row 7 (new): Resolving-Mac | Gdf7 | Bmpr1a | FAPs
# ---------------------------------------------------------------------
$ws.Cells.Item(7, 1).Value  = "Resolving-Mac"
$ws.Cells.Item(7, 2).Value  = "Gdf7"
$ws.Cells.Item(7, 3).Value  = "Bmpr1a"
$ws.Cells.Item(7, 4).Value  = "FAPs"
$ws.Cells.Item(7, 5).Value  = 1
$ws.Cells.Item(7, 6).Value  = 0.3333333333333333
$ws.Cells.Item(7, 7).Value  = 0.06826433333333333
$ws.Cells.Item(7, 8).Value  = 0.204793
$ws.Cells.Item(7, 9).Value  = 0.7008538526034804
$ws.Cells.Item(7, 10).Value = 0.7008538526034804
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 41.286995
$ws.Cells.Item(7, 14).Value = 123.860985
$ws.Cells.Item(7, 15).Value = 0.5401386314560596
$ws.Cells.Item(7, 16).Value = 0.5401386314560597
$ws.Cells.Item(7, 17).Value = 2.818429189011666
$ws.Cells.Item(7, 18).Value = 25.365862701105
$ws.Cells.Item(7, 19).Value = 0.3785582407959508
$ws.Cells.Item(7, 20).Value = 0.3785582407959509

# ---------------------------------------------------------------------
# Row 8 (new): Resolving-Mac | Gdf7 | Bmpr1a | MuSCs
# ---------------------------------------------------------------------
$ws.Cells.Item(8, 1).Value  = "Resolving-Mac"
$ws.Cells.Item(8, 2).Value  = "Gdf7"
$ws.Cells.Item(8, 3).Value  = "Bmpr1a"
$ws.Cells.Item(8, 4).Value  = "MuSCs"
$ws.Cells.Item(8, 5).Value  = 1
$ws.Cells.Item(8, 6).Value  = 0.3333333333333333
$ws.Cells.Item(8, 7).Value  = 0.06826433333333333
$ws.Cells.Item(8, 8).Value  = 0.204793
$ws.Cells.Item(8, 9).Value  = 0.7008538526034804
$ws.Cells.Item(8, 10).Value = 0.7008538526034804
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 27.73243066666667
$ws.Cells.Item(8, 14).Value = 83.197292
$ws.Cells.Item(8, 15).Value = 0.3628105447549136
$ws.Cells.Item(8, 16).Value = 0.3628105447549136
$ws.Cells.Item(8, 17).Value = 1.893135891172889
$ws.Cells.Item(8, 18).Value = 17.038223020556
$ws.Cells.Item(8, 19).Value = 0.2542771680566486
$ws.Cells.Item(8, 20).Value = 0.2542771680566486

# ---------------------------------------------------------------------
# Row 9 (new): Resolving-Mac | Gdf7 | Bmpr1a | Resolving-Mac
# ---------------------------------------------------------------------
$ws.Cells.Item(9, 1).Value  = "Resolving-Mac"
$ws.Cells.Item(9, 2).Value  = "Gdf7"
$ws.Cells.Item(9, 3).Value  = "Bmpr1a"
$ws.Cells.Item(9, 4).Value  = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value  = 1
$ws.Cells.Item(9, 6).Value  = 0.3333333333333333
$ws.Cells.Item(9, 7).Value  = 0.06826433333333333
$ws.Cells.Item(9, 8).Value  = 0.204793
$ws.Cells.Item(9, 9).Value  = 0.7008538526034804
$ws.Cells.Item(9, 10).Value = 0.7008538526034804
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.096171666666667
$ws.Cells.Item(9, 14).Value = 3.288515
$ws.Cells.Item(9, 15).Value = 0.01434070616847367
$ws.Cells.Item(9, 16).Value = 0.01434070616847367
$ws.Cells.Item(9, 17).Value = 0.0748294280438889
$ws.Cells.Item(9, 18).Value = 0.6734648523950001
$ws.Cells.Item(9, 19).Value = 0.01005073916722926
$ws.Cells.Item(9, 20).Value = 0.01005073916722926
